# Generate Report for Handoff
# Adds a new row (row 3) to each of the three worksheets (Overview, zh-cn, de-de)
# describing the newly-handed-off file "72de58a6-555a-418c-9d9d-a9ab13e66273.md".

$wb = $excel.ActiveWorkbook

$guidFile   = "72de58a6-555a-418c-9d9d-a9ab13e66273.md"
$repoBase   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/708156872113245be46da602ea16c0258f738bbd/e2e/"
$targetUrl  = $repoBase + $guidFile

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Range("A3").Value = $guidFile
$ws.Range("C3").Value = ".md"
$ws.Range("D3").Value = "'"
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-21 18:44:43"
$ws.Hyperlinks.Add($ws.Range("B3"), $targetUrl, [Type]::Missing, [Type]::Missing, ("e2e\" + $guidFile)) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "72de58a6-555a-418c-9d9d-a9ab13e66273.6b5f57d547faadb1efec930759edb568c4599aac.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-21 18:44:39"
$ws.Range("I3").Value = "'"
$ws.Range("J3").Value = "'"
$ws.Range("K3").Value = "0001-01-01 00:00:00"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"
$ws.Hyperlinks.Add($ws.Range("A3"), $targetUrl, [Type]::Missing, [Type]::Missing, $guidFile) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "72de58a6-555a-418c-9d9d-a9ab13e66273.6b5f57d547faadb1efec930759edb568c4599aac.de-de.xlf"
$ws.Range("H3").Value = "2016-08-21 18:44:43"
$ws.Range("I3").Value = "'"
$ws.Range("J3").Value = "'"
$ws.Range("K3").Value = "0001-01-01 00:00:00"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"
$ws.Hyperlinks.Add($ws.Range("A3"), $targetUrl, [Type]::Missing, [Type]::Missing, $guidFile) | Out-Null
